$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "XNPV" example row (row 7): label in G7, formula in H7
$ws.Range("G7").Value = "XNPV"
$ws.Range("H7").Formula = "=XNPV(0.05, B11:F11, B10:F10)"

# Update the selected cell to match the author's saved cursor position
[void]$ws.Range("H7").Select()
